$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 199 and 200 had their match data (columns F:V) swapped - the
#    "America De Cali - Bucaramanga" match and the "Atl. Nacional - Deportes
#    Tolima" match traded places in the sheet. Columns A:E (index / country /
#    tournament / season / match date) stay put.
# ---------------------------------------------------------------------------
$swapCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
foreach ($col in $swapCols) {
    $top = $ws.Range($col + "199").Value()
    $bottom = $ws.Range($col + "200").Value()
    $ws.Range($col + "199").Value = $bottom
    $ws.Range($col + "200").Value = $top
}

# ---------------------------------------------------------------------------
# 2) Two brand-new matches were appended at the bottom of the sheet (rows
#    214 and 215). Copy the formatting of the last existing row (213) down
#    so the new rows pick up the same styles (bold/boxed index column,
#    date-time formatted match-date column) already used throughout the
#    sheet, then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A213:V213").Copy()
$ws.Range("A214:V215").PasteSpecial(-4122)

function Set-TextCell($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

# --- row 214: Aguilas 1-1 Junior -------------------------------------------
$ws.Range("A214").Value = 213
$ws.Range("B214").Value = "colombia"
$ws.Range("C214").Value = "primera-a"
Set-TextCell "D214" "2023"
$ws.Range("E214").Value = 45260.01041666666
$ws.Range("F214").Value = "Aguilas"
$ws.Range("G214").Value = 1
$ws.Range("H214").Value = "Junior"
$ws.Range("I214").Value = 1
$ws.Range("J214").Value = 2.11
$ws.Range("K214").Value = "26/11/2023 01:42"
$ws.Range("L214").Value = 2.35
$ws.Range("M214").Value = "30/11/2023 00:11"
$ws.Range("N214").Value = 3.2
$ws.Range("O214").Value = "26/11/2023 01:42"
$ws.Range("P214").Value = 3.36
$ws.Range("Q214").Value = "30/11/2023 00:11"
$ws.Range("R214").Value = 3.9
$ws.Range("S214").Value = "26/11/2023 01:42"
$ws.Range("T214").Value = 3.23
$ws.Range("U214").Value = "30/11/2023 00:11"
$ws.Range("V214").Value = "https://www.betexplorer.com/football/colombia/primera-a/aguilas-doradas-junior/tEEOkEab/"

# --- row 215: Dep. Cali 0-2 Deportes Tolima ---------------------------------
$ws.Range("A215").Value = 214
$ws.Range("B215").Value = "colombia"
$ws.Range("C215").Value = "primera-a"
Set-TextCell "D215" "2023"
$ws.Range("E215").Value = 45260.10416666666
$ws.Range("F215").Value = "Dep. Cali"
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = "Deportes Tolima"
$ws.Range("I215").Value = 2
$ws.Range("J215").Value = 2.19
$ws.Range("K215").Value = "25/11/2023 23:12"
$ws.Range("L215").Value = 3.39
$ws.Range("M215").Value = "30/11/2023 02:26"
$ws.Range("N215").Value = 3.2
$ws.Range("O215").Value = "25/11/2023 23:12"
$ws.Range("P215").Value = 3.47
$ws.Range("Q215").Value = "30/11/2023 02:21"
$ws.Range("R215").Value = 3.66
$ws.Range("S215").Value = "25/11/2023 23:12"
$ws.Range("T215").Value = 2.22
$ws.Range("U215").Value = "30/11/2023 02:26"
$ws.Range("V215").Value = "https://www.betexplorer.com/football/colombia/primera-a/dep-cali-deportes-tolima/K4DSlYE4/"
